# Auto-generated: apply the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.513.08"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.839.13"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "261.01"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "0.5253"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "0.3193"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "0.06793"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").Value = "0.7842"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("D12").Value = "0.07756"
$ws.Range("E12").Value = "  +0.90%  "
$ws.Range("D13").Value = "1.839.61"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").Value = "87.87"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "5.013"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D17").Value = "13.85"
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "0.000007952"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "26.543.64"
$ws.Range("E20").Value = "  +0.20%  "
$ws.Range("D21").Value = "2.070.18"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").Value = "4.623"
$ws.Range("E22").Value = "  +1.12%  "
$ws.Range("D23").Value = "5.974"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "9.334"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("D25").Value = "141.87"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").Value = "2.196"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").Value = "1.679"
$ws.Range("E27").Value = "  +1.35%  "
$ws.Range("D28").Value = "16.95"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "111.56"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "4.164"
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("D31").Value = "0.08706"
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("D32").Value = "4.076"
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04880"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("D34").Value = "0.7273"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("D35").Value = "1.136"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "2.861"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").Value = "3.093"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").Value = "2.239"
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("D39").Value = "0.01754"
$ws.Range("E39").Value = "  -0.64%  "
$ws.Range("D40").Value = "0.4806"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("D41").Value = "0.8951"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "109.51"
$ws.Range("D43").Value = "5.937"
$ws.Range("E43").Value = "  -2.53%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "7.669"
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "0.4173"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("D47").Value = "8.963"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "0.05844"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "0.1233"
$ws.Range("D50").Value = "34.92"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("D51").Value = "0.8915"
$ws.Range("E51").Value = "  +0.96%  "
